$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Preserve numeric-looking Price strings as text (source cells are inlineStr, not numbers)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.759.37"
$ws.Range("E2").Value = "  +4.67%  "
$ws.Range("D3").Value = "2.632.78"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "567.92"
$ws.Range("E5").Value = "  +6.29%  "
$ws.Range("D6").Value = "146.33"
$ws.Range("E6").Value = "  +3.52%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "0.607"
$ws.Range("E8").Value = "  +3.15%  "
$ws.Range("D9").Value = "2.657.64"
$ws.Range("E9").Value = "  +3.09%  "
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("E11").Value = "  +5.87%  "
$ws.Range("E12").Value = "  +7.02%  "
$ws.Range("E13").Value = "  +4.22%  "
$ws.Range("D14").Value = "3.112.04"
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("D15").Value = "60.643.47"
$ws.Range("E15").Value = "  +4.53%  "
$ws.Range("D16").Value = "22.09"
$ws.Range("E16").Value = "  +6.87%  "
$ws.Range("E17").Value = "  +5.59%  "
$ws.Range("D18").Value = "2.654.58"
$ws.Range("E18").Value = "  +3.58%  "
$ws.Range("D19").Value = "4.54"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").Value = "343.50"
$ws.Range("E20").Value = "  +2.88%  "
$ws.Range("D21").Value = "10.45"
$ws.Range("E21").Value = "  +4.19%  "
$ws.Range("E22").Value = "  +3.73%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "66.12"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("E25").Value = "  +4.39%  "
$ws.Range("E26").Value = "  +3.28%  "
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("E28").Value = "  +5.58%  "
$ws.Range("D29").Value = "0.0₃0804"
$ws.Range("E29").Value = "  +11.16%  "
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Value = "1.71"
$ws.Range("E31").Value = "  +5.10%  "
$ws.Range("D32").Value = "6.15"
$ws.Range("E32").Value = "  +4.92%  "
$ws.Range("D33").Value = "158.97"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("D35").Value = "4.11"
$ws.Range("E35").Value = "  +6.30%  "
$ws.Range("D36").Value = "0.894"
$ws.Range("E36").Value = "  +7.20%  "
$ws.Range("E37").Value = "  +5.96%  "
$ws.Range("D38").Value = "0.890"
$ws.Range("E38").Value = "  +8.29%  "
$ws.Range("E39").Value = "  +7.86%  "
$ws.Range("D40").Value = "37.46"
$ws.Range("E40").Value = "  +1.49%  "
$ws.Range("D41").Value = "297.93"
$ws.Range("E41").Value = "  +6.07%  "
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").Value = "0.0984"
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("D45").Value = "0.603"
$ws.Range("E45").Value = "  +2.54%  "
$ws.Range("E46").Value = "  +2.63%  "
$ws.Range("D47").Value = "19.52"
$ws.Range("E47").Value = "  +3.88%  "
$ws.Range("D48").Value = "128.07"
$ws.Range("E48").Value = "  +16.16%  "
$ws.Range("D49").Value = "10.70"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("E50").Value = "  +4.33%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "4.71"
$ws.Range("E51").Value = "  +7.05%  "
